$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "91.500.30"
$ws.Range("E2").Value = "  +2.55%  "
$ws.Range("D3").Value = "3.131.96"
$ws.Range("E3").Value = "  +1.07%  "
$ws.Range("E4").Value = "  +0.17%  "
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "218.49"
$r.ClearFormats()
$ws.Range("E5").Value = "  +2.43%  "
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = "625.10"
$r.ClearFormats()
$ws.Range("E6").Value = "  +0.28%  "
$r = $ws.Range("D7")
$r.NumberFormat = "@"
$r.Value = "1.07"
$r.ClearFormats()
$ws.Range("E7").Value = "  +30.52%  "
$r = $ws.Range("D8")
$r.NumberFormat = "@"
$r.Value = "0.377"
$r.ClearFormats()
$ws.Range("E8").Value = "  +0.71%  "
$r = $ws.Range("D9")
$r.NumberFormat = "@"
$r.Value = "1.00"
$r.ClearFormats()
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("D10").Value = "3.127.96"
$ws.Range("E10").Value = "  +1.01%  "
$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = "0.771"
$r.ClearFormats()
$ws.Range("E11").Value = "  +24.17%  "
$r = $ws.Range("D12")
$r.NumberFormat = "@"
$r.Value = "0.195"
$r.ClearFormats()
$ws.Range("E12").Value = "  +7.61%  "
$r = $ws.Range("D13")
$r.NumberFormat = "@"
$r.Value = "0.0000252"
$r.ClearFormats()
$ws.Range("E13").Value = "  +4.30%  "
$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = "34.96"
$r.ClearFormats()
$ws.Range("E14").Value = "  +8.24%  "
$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = "5.56"
$r.ClearFormats()
$ws.Range("E15").Value = "  +4.69%  "
$ws.Range("D16").Value = "91.201.60"
$ws.Range("E16").Value = "  +2.58%  "
$ws.Range("D17").Value = "3.705.91"
$ws.Range("E17").Value = "  +1.13%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.124.44"
$ws.Range("E18").Value = "  +0.69%  "
$ws.Range("B19").Value = "SuiNetwork"
$ws.Range("C19").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = "3.86"
$r.ClearFormats()
$ws.Range("E19").Value = "  +13.77%  "
$ws.Range("E20").Value = "  +2.89%  "
$ws.Range("E21").Value = "  +5.37%  "
$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = "446.65"
$r.ClearFormats()
$ws.Range("E22").Value = "  +5.36%  "
$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = "8.86"
$r.ClearFormats()
$ws.Range("E23").Value = "  +6.77%  "
$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = "5.23"
$r.ClearFormats()
$ws.Range("E24").Value = "  +5.72%  "
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = "6.28"
$r.ClearFormats()
$ws.Range("E25").Value = "  +12.87%  "
$r = $ws.Range("D26")
$r.NumberFormat = "@"
$r.Value = "90.16"
$r.ClearFormats()
$ws.Range("E26").Value = "  +9.34%  "
$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = "12.53"
$r.ClearFormats()
$ws.Range("E27").Value = "  +4.98%  "
$ws.Range("D28").Value = "3.292.06"
$ws.Range("E28").Value = "  +1.74%  "
$ws.Range("E29").Value = "  +0.23%  "
$ws.Range("E31").Value = "  +13.24%  "
$ws.Range("B32").Value = "Binance-PegBSC-USD"
$ws.Range("C32").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$r = $ws.Range("D32")
$r.NumberFormat = "@"
$r.Value = "0.908"
$r.ClearFormats()
$ws.Range("E32").Value = "  -15.76%  "
$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$r = $ws.Range("D33")
$r.NumberFormat = "@"
$r.Value = "527.21"
$r.ClearFormats()
$ws.Range("E33").Value = "  +3.38%  "
$r = $ws.Range("D34")
$r.NumberFormat = "@"
$r.Value = "3.75"
$r.ClearFormats()
$ws.Range("E34").Value = "  +1.84%  "
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$r = $ws.Range("D35")
$r.NumberFormat = "@"
$r.Value = "0.148"
$r.ClearFormats()
$ws.Range("E35").Value = "  +12.40%  "
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$r = $ws.Range("D36")
$r.NumberFormat = "@"
$r.Value = "7.13"
$r.ClearFormats()
$ws.Range("E36").Value = "  +6.04%  "
$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = "24.66"
$r.ClearFormats()
$ws.Range("E37").Value = "  +10.42%  "
$ws.Range("E38").Value = "  +4.65%  "
$ws.Range("E39").Value = "  +4.16%  "
$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = "0.164"
$r.ClearFormats()
$ws.Range("E40").Value = "  +24.38%  "
$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = "0.0869"
$r.ClearFormats()
$ws.Range("E41").Value = "  +24.99%  "
$ws.Range("B42").Value = "WhiteBITCoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = "22.30"
$r.ClearFormats()
$ws.Range("E42").Value = "  +0.09%  "
$ws.Range("E43").Value = "  -0.11%  "
$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = "0.409"
$r.ClearFormats()
$ws.Range("E44").Value = "  +12.15%  "
$ws.Range("E45").Value = "  +6.57%  "
$ws.Range("E46").Value = "  -0.02%  "
$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = "149.54"
$r.ClearFormats()
$ws.Range("E47").Value = "  +2.71%  "
$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = "44.39"
$r.ClearFormats()
$ws.Range("E48").Value = "  +2.46%  "
$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = "1.32"
$r.ClearFormats()
$ws.Range("E49").Value = "  +8.90%  "
$ws.Range("E50").Value = "  +9.39%  "
$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = "169.48"
$r.ClearFormats()
$ws.Range("E51").Value = "  +5.20%  "
